# Edit script: apply corrected classification values to Sheet1, row 1.
# The original commit message indicates "error in code has been corrected",
# i.e. a set of misclassified cell values (1 <-> 2) are being fixed in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H1").Value = 2
$ws.Range("W1").Value = 2
$ws.Range("X1").Value = 1
$ws.Range("BA1").Value = 1
$ws.Range("BY1").Value = 1
$ws.Range("CB1").Value = 1
$ws.Range("CJ1").Value = 1
$ws.Range("CQ1").Value = 1
$ws.Range("EA1").Value = 2
$ws.Range("FB1").Value = 1
$ws.Range("FK1").Value = 2
$ws.Range("FZ1").Value = 2
$ws.Range("GJ1").Value = 2
$ws.Range("GR1").Value = 2
$ws.Range("HX1").Value = 1
$ws.Range("KB1").Value = 2
$ws.Range("KM1").Value = 1
$ws.Range("KV1").Value = 1
$ws.Range("MB1").Value = 2
$ws.Range("MG1").Value = 1
$ws.Range("MI1").Value = 2
$ws.Range("NJ1").Value = 1
$ws.Range("NW1").Value = 2
$ws.Range("NY1").Value = 1
$ws.Range("OB1").Value = 1
$ws.Range("OD1").Value = 2
$ws.Range("OF1").Value = 2
$ws.Range("OL1").Value = 1
$ws.Range("OW1").Value = 2
$ws.Range("OX1").Value = 2
$ws.Range("OY1").Value = 1
$ws.Range("OZ1").Value = 1
$ws.Range("PD1").Value = 1
$ws.Range("PK1").Value = 2
$ws.Range("PO1").Value = 1
$ws.Range("QI1").Value = 2
$ws.Range("SN1").Value = 1
$ws.Range("TV1").Value = 2
$ws.Range("UD1").Value = 1
$ws.Range("UI1").Value = 2
$ws.Range("UT1").Value = 1
$ws.Range("VB1").Value = 2
$ws.Range("VG1").Value = 1
$ws.Range("VI1").Value = 1
$ws.Range("WA1").Value = 2
$ws.Range("WL1").Value = 2
$ws.Range("WT1").Value = 2
$ws.Range("XX1").Value = 1
$ws.Range("XZ1").Value = 1
$ws.Range("YT1").Value = 2
$ws.Range("ZS1").Value = 2
$ws.Range("AAB1").Value = 2
$ws.Range("ACY1").Value = 2
$ws.Range("ADB1").Value = 1
$ws.Range("AFX1").Value = 1
$ws.Range("AGK1").Value = 2
$ws.Range("AHD1").Value = 2
$ws.Range("AHH1").Value = 1
$ws.Range("AIE1").Value = 2
$ws.Range("AIJ1").Value = 2
$ws.Range("AJD1").Value = 1
$ws.Range("AJN1").Value = 2
$ws.Range("AJW1").Value = 2
$ws.Range("AKE1").Value = 1
$ws.Range("AKG1").Value = 1

$wb.Save()
